$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "no" -> "el" text used in B2 (ablation central-range label)
$ws.Range("B2").Value = "el"

# Update the ablation importance values in B3:B136
$ws.Range("B3").Value = 0.02644517458975315
$ws.Range("B4").Value = 0.02858145907521248
$ws.Range("B5").Value = 0.05854104086756706
$ws.Range("B6").Value = 0.0438130758702755
$ws.Range("B7").Value = 0.04786625504493713
$ws.Range("B8").Value = 0.005768594797700644
$ws.Range("B9").Value = 0.02866285853087902
$ws.Range("B10").Value = 0.02605069428682327
$ws.Range("B11").Value = -0.02363945171236992
$ws.Range("B12").Value = -0.00112161785364151
$ws.Range("B13").Value = 0.02742592431604862
$ws.Range("B14").Value = 0.004196058958768845
$ws.Range("B15").Value = 0.03863755986094475
$ws.Range("B16").Value = 0.005860290490090847
$ws.Range("B17").Value = -0.02680091559886932
$ws.Range("B18").Value = -0.02547126635909081
$ws.Range("B19").Value = -0.08242511749267578
$ws.Range("B20").Value = -0.05430565774440765
$ws.Range("B21").Value = -0.02056302689015865
$ws.Range("B22").Value = 0.01150072738528252
$ws.Range("B23").Value = 0.02650769799947739
$ws.Range("B24").Value = -0.07159511744976044
$ws.Range("B25").Value = 0.0579531155526638
$ws.Range("B26").Value = -0.004990175366401672
$ws.Range("B27").Value = -0.02592606842517853
$ws.Range("B28").Value = -0.02624731510877609
$ws.Range("B29").Value = 0.03094606846570969
$ws.Range("B30").Value = 0.03860154002904892
$ws.Range("B31").Value = 0.0159921683371067
$ws.Range("B32").Value = 0.08934935927391052
$ws.Range("B33").Value = 0.04135732352733612
$ws.Range("B34").Value = -0.03367913886904716
$ws.Range("B35").Value = 0.01026834733784199
$ws.Range("B36").Value = 0.02395711094141006
$ws.Range("B37").Value = -0.03062940761446953
$ws.Range("B38").Value = -0.008171543478965759
$ws.Range("B39").Value = -0.0242218729108572
$ws.Range("B40").Value = 0.02263829857110977
$ws.Range("B41").Value = 0.003193280193954706
$ws.Range("B42").Value = 0.08152377605438232
$ws.Range("B43").Value = 0.05306227505207062
$ws.Range("B44").Value = -0.06043200194835663
$ws.Range("B45").Value = -0.01713793352246284
$ws.Range("B46").Value = 0.06619825214147568
$ws.Range("B47").Value = 0.03948301821947098
$ws.Range("B48").Value = -0.007296753581613302
$ws.Range("B49").Value = 0.04829274863004684
$ws.Range("B50").Value = 0.01767875626683235
$ws.Range("B51").Value = 0.0003033801913261414
$ws.Range("B52").Value = -0.03591246530413628
$ws.Range("B53").Value = 0.02717983722686768
$ws.Range("B54").Value = -0.04839520156383514
$ws.Range("B55").Value = -0.00758073478937149
$ws.Range("B56").Value = -0.06363990157842636
$ws.Range("B57").Value = -0.01565032266080379
$ws.Range("B58").Value = -0.004895614460110664
$ws.Range("B59").Value = -0.01080658286809921
$ws.Range("B60").Value = -0.02818277105689049
$ws.Range("B61").Value = -0.04547712206840515
$ws.Range("B62").Value = -0.02017189189791679
$ws.Range("B63").Value = -0.05626610666513443
$ws.Range("B64").Value = -0.00492263026535511
$ws.Range("B65").Value = -0.05653904005885124
$ws.Range("B66").Value = -0.03268317878246307
$ws.Range("B67").Value = 0.05089400336146355
$ws.Range("B68").Value = -0.05173088610172272
$ws.Range("B69").Value = -0.0129594411700964
$ws.Range("B70").Value = 0.05494849756360054
$ws.Range("B71").Value = -0.03833058848977089
$ws.Range("B72").Value = 0.02440397255122662
$ws.Range("B73").Value = 0.05478239059448242
$ws.Range("B74").Value = -0.03656921908259392
$ws.Range("B75").Value = -0.04076679795980453
$ws.Range("B76").Value = 0.04844409599900246
$ws.Range("B77").Value = -0.02951251901686192
$ws.Range("B78").Value = -0.003232178278267384
$ws.Range("B79").Value = 0.01473437249660492
$ws.Range("B80").Value = 0.005027965176850557
$ws.Range("B81").Value = -0.0009829029440879822
$ws.Range("B82").Value = 0.005836565047502518
$ws.Range("B83").Value = -0.02832100726664066
$ws.Range("B84").Value = -0.01563811302185059
$ws.Range("B85").Value = -0.0003687338903546333
$ws.Range("B86").Value = -0.03015160746872425
$ws.Range("B87").Value = 0.04032589867711067
$ws.Range("B88").Value = -0.04026233032345772
$ws.Range("B89").Value = 0.03141788765788078
$ws.Range("B90").Value = -0.08803508430719376
$ws.Range("B91").Value = -0.03992601856589317
$ws.Range("B92").Value = 0.03762850537896156
$ws.Range("B93").Value = -0.02673275396227837
$ws.Range("B94").Value = -0.04942956939339638
$ws.Range("B95").Value = -0.03061344474554062
$ws.Range("B96").Value = -0.1117067113518715
$ws.Range("B97").Value = 0.03908052295446396
$ws.Range("B98").Value = 0.02212160266935825
$ws.Range("B99").Value = 0.0405014231801033
$ws.Range("B100").Value = -0.04356152564287186
$ws.Range("B101").Value = -0.03478874266147614
$ws.Range("B102").Value = 0.01153452601283789
$ws.Range("B103").Value = -0.04361241310834885
$ws.Range("B104").Value = -0.02877399884164333
$ws.Range("B105").Value = 0.01876598037779331
$ws.Range("B106").Value = -0.02730368264019489
$ws.Range("B107").Value = -0.0356324128806591
$ws.Range("B108").Value = 0.00756595004349947
$ws.Range("B109").Value = 0.002646961947903037
$ws.Range("B110").Value = -0.04497546330094337
$ws.Range("B111").Value = -0.101141981780529
$ws.Range("B112").Value = -0.06047642603516579
$ws.Range("B113").Value = 0.01402320060878992
$ws.Range("B114").Value = 0.03536102175712585
$ws.Range("B115").Value = -0.1131041795015335
$ws.Range("B116").Value = -0.01362222526222467
$ws.Range("B117").Value = 0.04522600397467613
$ws.Range("B118").Value = 0.007307623978704214
$ws.Range("B119").Value = -0.02908360213041306
$ws.Range("B120").Value = 0.01935539022088051
$ws.Range("B121").Value = -0.0524100661277771
$ws.Range("B122").Value = -0.06680228561162949
$ws.Range("B123").Value = 0.05929967015981674
$ws.Range("B124").Value = -0.002190083265304565
$ws.Range("B125").Value = -0.02165385708212852
$ws.Range("B126").Value = -0.07507848739624023
$ws.Range("B127").Value = -0.01523554883897305
$ws.Range("B128").Value = -0.01213733106851578
$ws.Range("B129").Value = -0.04678767174482346
$ws.Range("B130").Value = -0.01741739548742771
$ws.Range("B131").Value = -0.03670401871204376
$ws.Range("B132").Value = -0.019779397174716
$ws.Range("B133").Value = 0.03499685227870941
$ws.Range("B134").Value = -0.03210711106657982
$ws.Range("B135").Value = -0.03876043483614922
$ws.Range("B136").Value = 0.04722781106829643
